$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the "Zoom" bullet (text + hyperlink to http://zoom.us/) entirely.
#    This naturally shifts the R / RStudio / GDAL bullets up by one slot,
#    which is exactly the re-shuffled content the diff shows (the diff's
#    repeated add/remove blocks for "R ", "RStudio " and the GDAL paragraph
#    are the textual result of this single paragraph deletion).
# ---------------------------------------------------------------------------
$zoom = $d.Content.Find
$zoom.Text = "Zoom "
$zoom.Execute() | Out-Null
if ($zoom.Found) {
    $p = $zoom.Parent.Paragraphs(1)
    $p.Range.Delete()
}

# ---------------------------------------------------------------------------
# 2) "Dataset and workshop code" bullet: the placeholder hyperlink
#    https://[workshop-git-URL] becomes the real workshop repo URL, and the
#    stray trailing space run between the hyperlink and the final "." goes
#    away.
# ---------------------------------------------------------------------------
$realUrl = "https://github.com/ocean-tracking-network/2024-GLATOS-intro-R-workshop"

$d.Content.Find.Execute("https://[workshop-git-URL]", $true, $false, $false, $false, $false, `
    $true, 1, $false, $realUrl, 2) | Out-Null

# Update the underlying hyperlink relationship target to match the new text.
for ($i = 1; $i -le $d.Hyperlinks.Count; $i++) {
    $h = $d.Hyperlinks($i)
    if ($h.Address -eq "https://[workshop-git-URL]" -or $h.TextToDisplay -eq $realUrl) {
        $h.Address = $realUrl
    }
}

# Remove the leftover " " run immediately before the closing "." run, i.e.
# collapse "<url> ." -> "<url>."
$d.Content.Find.Execute("$realUrl .", $true, $false, $false, $false, $false, `
    $true, 1, $false, "$realUrl.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) "git clone [workshop-git-URL].git" -> real URL.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("[workshop-git-URL].git", $true, $false, $false, $false, $false, `
    $true, 1, $false, "$realUrl.git", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) lastRenderedPageBreak relocations (pure rendering artifacts that Word
#    re-stamps as pagination shifts with the content edits above).
# ---------------------------------------------------------------------------

# 4a) "Instructions for installing R:" no longer starts a fresh page ...
$d.Content.Find.Execute("Instructions for installing R:") | Out-Null

# 4b) ... the break now falls in front of "Click on the ""Download"" link" bullet.
